$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data layout: row number, then values for columns C,D,E,F,G,H,I,J,K,M,N (in order)
$data = @"
2|4.409979911005498,8.103642006615821,12.743477862311,36.3220944133288,40.95254160811817,17.19658491585501,25.47273313125103,9.88317813655183,15.52051984591185,17.88226560424409,18.67081359058918
3|4.354570677511079,8.09585827525177,12.76034808263853,36.33872941042381,40.89799378580049,17.24063348324793,25.49418085073745,9.907994448905029,15.14732228422118,17.74553255921584,18.73461801456704
4|4.319705230154812,8.091926194027732,12.77254663999969,36.35887737466083,40.87867573867039,17.27125000468269,25.51420518833381,9.924322079756019,14.9163632322336,17.66426631682183,18.77569587930289
5|4.305292462703618,8.090538212922299,12.77798051683051,36.36958085464674,40.87436640163381,17.28462213632875,25.52408489296537,9.931250221743388,14.82192997036291,17.63185479672971,18.79291506630276
6|4.30288708333496,8.090320725199726,12.77891076668978,36.37150855762013,40.87386591366981,17.28689659758314,25.52582913109169,9.932417223366672,14.80623455752213,17.62651630396103,18.79580330943132
7|4.319511673909799,8.091906605506979,12.77261804889903,36.3590116391263,40.87860320016139,17.2714267224847,25.51433147306949,9.924414403223288,14.9150907600936,17.66382631112056,18.77592615942275
8|4.391054349618026,8.100783178167021,12.7489128867226,36.32576580430397,40.93079194877622,17.21103051713067,25.47870320580893,9.89150870522837,15.39230045060593,17.83457989219173,18.69241952051017
9|4.524339190559011,8.124850429442056,12.7170205679736,36.33955187263147,41.14548397002953,17.12101826132984,25.46338111690952,9.835617082958883,16.30784433114446,18.18942735751998,18.54368677654554
10|4.617587249491145,8.146507080335541,12.70247353229891,36.3979542577019,41.37128906007108,17.07233749746541,25.48552633424875,9.799798186896187,16.96077947526372,18.46048051385055,18.44347930579238
11|4.658912216636433,8.157202359941234,12.69778085618261,36.4349985712805,41.48861957897451,17.05400448998229,25.50286567347451,9.784637869646408,17.25213288358097,18.58562892201046,18.39984154399553
12|4.67439712076397,8.1613718402659,12.69628016317885,36.45052935978433,41.53512935787221,17.04761214020127,25.51047555751503,9.779059796331731,17.36153970257039,18.63324974920228,18.38359562018288
13|4.671069564750969,8.160468588282805,12.69659108311462,36.44711776023797,41.52502055733569,17.04896436061311,25.50879023274165,9.78025389698659,17.33801962564559,18.62298409635648,18.38708209267397
14|4.660189492162214,8.157543003570533,12.69765185834081,36.43624620491225,41.492404433605,17.05346755578935,25.50347083175871,9.784175696605153,17.26115291793604,18.58954231355838,18.39849940256308
15|4.65350358818233,8.155766491617202,12.69833758373981,36.42978264972352,41.47269618413774,17.05629756449085,25.50034844078059,9.786599107707968,17.21394668315859,18.56908713226288,18.40552909938753
16|4.614863976694927,8.145824933295764,12.70281890428081,36.39574387221923,41.36391347110471,17.07361246822278,25.48453931272721,9.800811745128556,16.94161540678049,18.45233587007007,18.44637021036643
17|4.59087480310801,8.139940844381412,12.70606074005911,36.37754349168588,41.30090674013906,17.08521229388825,25.47670145547005,9.809820978420895,16.77301792682126,18.38116013740328,18.47192272149574
18|4.576974268519213,8.136635990458466,12.70810656673999,36.3680615634414,41.26604371590712,17.0922429239606,25.47287729857053,9.815109591131556,16.67551729058484,18.3403975375408,18.48680319510533
19|4.572250372205943,8.135530737401938,12.70883038418866,36.36502065206404,41.25447675551261,17.09468492766408,25.47169998764495,9.816918564646135,16.6424179950175,18.3266272920595,18.49187299574509
20|4.593439158917556,8.140558999298731,12.70569688984824,36.37937888033272,41.30747156730173,17.0839403333069,25.47746502340017,9.808850884721902,16.79102084627252,18.38871897730308,18.46918364659087
21|4.663389734714441,8.158399093449303,12.69733278779014,36.43939869309188,41.50192834212284,17.05212991893568,25.50500495130776,9.783019352269825,17.28375637832999,18.59935902153211,18.3951383065163
22|4.708147887758302,8.170753704904973,12.6934767025855,36.48738233692319,41.64112374811181,17.03454629207201,25.5290870943465,9.767085824909241,17.60036735539369,18.73834761204813,18.348369599356
23|4.68434945271817,8.164096862591459,12.69538759745207,36.46097291773431,41.56573297536287,17.04363707732452,25.51567795755878,9.775503104515506,17.43191526916332,18.66405747350718,18.37318272951733
24|4.592280152160013,8.140279288412692,12.70586081942377,36.37854604279703,41.30449936971654,17.08451425997329,25.47711769003541,9.809289124268023,16.78288349884758,18.38530113315602,18.47042139039573
25|4.489075074261406,8.117635097336473,12.72408693605649,36.32734763629785,41.07540345267047,17.14231271425812,25.46166733753331,9.849814783288521,16.06312437495759,18.09147800995509,18.5823241666059
"@

$cols = @("C","D","E","F","G","H","I","J","K","M","N")

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split "\|"
    $row = [int]$parts[0]
    $vals = $parts[1] -split ","
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = [double]$vals[$i]
        $ws.Range("$col$row").Value = $val
    }
}
